$d = $word.ActiveDocument

# 1. Normal style: add <w:pPr><w:spacing w:after="120"/></w:pPr>
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceAfter = 6

# 2. Add new paragraph style "tei_collation" (styleId "teicollation"),
#    based on Heading4, followed by Normal, quick-styled.
$collation = $d.Styles.Add("teicollation", 1)
$collation.NameLocal = "tei_collation"
$collation.BaseStyle = $d.Styles("Heading4")
$collation.NextParagraphStyle = $d.Styles("Normal")
$collation.QuickStyle = $true

# 3. Add new paragraph style "tei_extent" (styleId "teiextent"),
#    based on Heading4, followed by Normal, quick-styled.
$extent = $d.Styles.Add("teiextent", 1)
$extent.NameLocal = "tei_extent"
$extent.BaseStyle = $d.Styles("Heading4")
$extent.NextParagraphStyle = $d.Styles("Normal")
$extent.QuickStyle = $true
